$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: becomes what was row 3 (Fecha 44350, Arica y Parinacota origin)
$ws.Range("D2").Value = 44350
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 167

# Row 3: becomes what was row 4 (Fecha 44291, Limache origin)
$ws.Range("D3").Value = 44291
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 11000
$ws.Range("M3").Value = 11000
$ws.Range("O3").Value = "Limache"
$ws.Range("P3").Value = 183

# Row 4: becomes what was row 2 (Fecha 44273, Provincia de Limarí origin)
$ws.Range("D4").Value = 44273
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 14000
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 233
